$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 and Row 17 hold the two "estado de cuenta" detail records.
# Update them so row 16 carries GERMAN ENRIQUE DUPERRET TRESPALACIOS's
# data and row 17 carries YENIS ESTHER CASTELLAR CASTELLAR's data
# (the two records swap places), and add the new worker's period values.

$ws.Range("C16").Value = "73190466"
$ws.Range("D16").Value = "GERMAN ENRIQUE DUPERRET TRESPALACIOS"
$ws.Range("E16").Value = "2001"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

$ws.Range("C17").Value = "45757793"
$ws.Range("D17").Value = "YENIS ESTHER CASTELLAR CASTELLAR"
$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

$wb.Save()
